$wb = $excel.ActiveWorkbook

# --- "data" sheet: insert a new "Flowbuds" column before "Flowers" ---
$wsData = $wb.Worksheets.Item("data")
$wsData.Columns("G:G").Insert()
$wsData.Range("G1").Value = "Flowbuds"
$wsData.Range("G2").Value = "N"

# --- "meta" sheet: reword the Flowbuds description, and add a trailing
#     formatted (but empty) cell at A25 matching the other row-label cells ---
$wsMeta = $wb.Worksheets.Item("meta")
$wsMeta.Range("B14").Value = "Does the plant have any flower bud (not open)?"
$wsMeta.Range("A16").Copy()
$wsMeta.Range("A25").PasteSpecial(-4122)

# --- cosmetic view state: zoom + selected cell on each sheet ---
$wsData.Activate()
$excel.ActiveWindow.Zoom = 125
$wsData.Range("F16").Select() | Out-Null

$wsMeta.Activate()
$excel.ActiveWindow.Zoom = 125
$wsMeta.Range("B15").Select() | Out-Null
